# Revision to experiments problem set
# Move "Songtao" on the "Lead discussion" sheet from week 5 (B6) to week 7 (B8),
# and make "Lead discussion" the active/selected sheet with B8 selected.

$wb = $excel.ActiveWorkbook

$wsLead = $wb.Worksheets.Item("Lead discussion")

# Clear the old entry and write the new one.
$wsLead.Range("B6").Value = ""
$wsLead.Range("B8").Value = "Songtao"

# Activate the "Lead discussion" sheet and select B8, matching the
# new tabSelected/activeTab + selection state in the workbook.
$wsLead.Activate()
$wsLead.Range("B8").Select()
